$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of trading data for 2021-11-12
$ws.Range("A6").Value = 44512
$ws.Range("B6").Value = 3504
$ws.Range("C6").Value = -0.04
$ws.Range("D6").Value = -258

# Match formatting of the row above (date style on column A)
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
